# Insert a new data row at row 22 (shifts existing rows 22-35 down to 23-36,
# carrying their values/formatting along automatically), then populate the
# newly inserted row with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("22:22").Insert()

$ws.Range("A22").Value = 11
$ws.Range("B22").Value = 'Vega Monumental Concepción'
$ws.Range("C22").Value = 'Bíobío'
$ws.Range("D22").Value = 44574
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112031
$ws.Range("G22").Value = 'Poroto verde'
$ws.Range("H22").Value = 'Magnum'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 27000
$ws.Range("L22").Value = 28000
$ws.Range("M22").Value = 27500
$ws.Range("N22").Value = '$/saco 25 kilos'
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 1100
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = 'Hortaliza'
